$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.269.34'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.610.38'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.89'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.486'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('E10').Value = '  +2.26%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.834.64'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.622.28'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.266.42'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.08'
$ws.Range('E17').Value = '  +3.15%  '
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '200.94'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.90'
$ws.Range('E24').Value = '  +4.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.16'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.23'
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0501'
$ws.Range('E30').Value = '  +6.29%  '
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('E32').Value = '  +2.76%  '
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E34').Value = '  +1.39%  '
$ws.Range('E35').Value = '  +0.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.157.07'
$ws.Range('E36').Value = '  +2.65%  '
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('E40').Value = '  +0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.495'
$ws.Range('E41').Value = '  +1.46%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.784'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.33'
$ws.Range('E43').Value = '  +3.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.745.75'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.63'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('E46').Value = '  +13.77%  '
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.79'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0507'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.407'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('E51').Value = '  -0.18%  '
